# Feat: Adiciona validador visual de inputs (Sanity Check)
#
# Updates the "Inatividade_Alunos" sheet (sheet11.xml) with refreshed data
# (new snapshot of inactive students), removes two rows that are no longer
# part of the report, widens column F, and introduces a new risk-level label
# "(4) Risco Inicial (10-15 dias)" used by the updated rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inatividade_Alunos")

# --- Row 9 / Row 10: swap the two "Laranja" students ---
$ws.Cells.Item(9, 1).Value = "Maria Luísa Pasinato Ferreira"
$ws.Cells.Item(9, 3).Value = "Boyan Sirakov"

$ws.Cells.Item(10, 1).Value = "Boris Leon Fontes Cardoso Bazan"
$ws.Cells.Item(10, 3).Value = "Paulo Mendes"

# --- Row 11: new student now flagged as "Risco Inicial" ---
$ws.Cells.Item(11, 1).Value = "Eduardo Sardenberg Tavares"
$ws.Cells.Item(11, 2).Value = "Aluno de mestrado"
$ws.Cells.Item(11, 3).Value = "Juliana Pereira"
$ws.Cells.Item(11, 4).Value = 45971
$ws.Cells.Item(11, 5).Value = 20
$ws.Cells.Item(11, 6).Value = "(4) Risco Inicial (10-15 dias)"

# --- Row 12: new student now flagged as "Risco Inicial" ---
$ws.Cells.Item(12, 1).Value = "Lucas Gomes Maddalena"
$ws.Cells.Item(12, 2).Value = "Aluno de mestrado"
$ws.Cells.Item(12, 3).Value = "Fernanda Baião"
$ws.Cells.Item(12, 4).Value = 45971
$ws.Cells.Item(12, 5).Value = 20
$ws.Cells.Item(12, 6).Value = "(4) Risco Inicial (10-15 dias)"

# --- Row 13 ---
$ws.Cells.Item(13, 1).Value = "Gabriel de Oliveira Esteves Dias "
$ws.Cells.Item(13, 2).Value = "Aluno de Doutorado"
$ws.Cells.Item(13, 3).Value = "Thiago Guerreiro"
$ws.Cells.Item(13, 4).Value = 45978
$ws.Cells.Item(13, 5).Value = 13
$ws.Cells.Item(13, 6).Value = "(1) Amarelo (> 15 dias)"

# --- Row 14 ---
$ws.Cells.Item(14, 1).Value = "Christopher Silva Aguiar"
$ws.Cells.Item(14, 2).Value = "Aluno de Doutorado"
$ws.Cells.Item(14, 3).Value = "Boyan Sirakov"
$ws.Cells.Item(14, 4).Value = 45978
$ws.Cells.Item(14, 5).Value = 13
$ws.Cells.Item(14, 6).Value = "(1) Amarelo (> 15 dias)"

# --- Row 15 ---
$ws.Cells.Item(15, 1).Value = "Victor Hugo dos Santos de Castro Marques"
$ws.Cells.Item(15, 2).Value = "Aluno de graduação"
$ws.Cells.Item(15, 3).Value = "Paulo Mendes"
$ws.Cells.Item(15, 4).Value = 45978
$ws.Cells.Item(15, 5).Value = 13
$ws.Cells.Item(15, 6).Value = "(1) Amarelo (> 15 dias)"

# --- Rows 16 and 17 no longer exist in the refreshed report ---
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(16).Delete()

# --- Column F (Nível de Risco) needs to be wider to fit the new label ---
$ws.Columns.Item(6).ColumnWidth = 31.8

# --- Shrink the conditional formatting range to match the new data extent ---
$oldRange = $ws.Range("F2:F18")
$fc = $oldRange.FormatConditions
$newRange = $ws.Range("F2:F16")
for ($i = 1; $i -le $fc.Count; $i++) {
    $fc.Item($i).ModifyAppliesToRange($newRange)
}
